# Remove DUNS-related columns from the CSV response spec:
#   Column C  -> "ueiDUNS"
#   Column X  -> "immediateParentEntity.ueiDUNS"
# (mirrors the commit "Removed duns information from the spec files")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reflect the author's selection (entire column W) prior to deleting,
# matching the saved sheet view state.
$ws.Range("W:W").Select()

# Delete the rightmost offending column first so column letters for the
# earlier column remain valid.
$ws.Range("X:X").Delete()
$ws.Range("C:C").Delete()
